# Update underlying (unrounded) simulation values on the hit_miss_rule
# sheet. Downstream cells (D/E/F columns) hold ROUND()/SUM() formulas over
# these raw values and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = 91.92653656005859
$ws.Range("I6").Value = 8.073463439941406

$ws.Range("H8").Value = 2.312006950378418
$ws.Range("I8").Value = 3.177914142608643

$ws.Range("H9").Value = 3.688869476318359
$ws.Range("I9").Value = 4.559158802032471

$ws.Range("H10").Value = 86.90134429931641
$ws.Range("I10").Value = 22.31196594238281
$ws.Range("J10").Value = 79.59637451171875

# Force a full recalculation on next load, mirroring the source workbook
# change (calcPr fullCalcOnLoad="1").
$wb.ForceFullCalculation = $true
